# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Santa Lucia" / "Timor Oriental" shared-string order ---
# Row 207 held "Santa Lucia", row 208 held "Timor Oriental".
# After the edit the shared-string table has "Timor Oriental" before
# "Santa Lucia", while the cells keep referencing the same shared-string
# slots, so the displayed text on rows 207/208 swaps.
$ws.Cells.Item(207, 1).Value = "Timor Oriental"
$ws.Cells.Item(208, 1).Value = "Santa Lucia"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 23:56"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 7317312
$ws.Cells.Item(4, 3).Value = 29751
$ws.Cells.Item(4, 4).Value = 4540875
$ws.Cells.Item(4, 5).Value = 2567014
$ws.Cells.Item(4, 7).Value = 246
$ws.Cells.Item(4, 8).Value = 209423

# --- Row 6: Brasil ---
$ws.Cells.Item(6, 2).Value = 4732309
$ws.Cells.Item(6, 3).Value = 14194
$ws.Cells.Item(6, 5).Value = 539731
$ws.Cells.Item(6, 7).Value = 300
$ws.Cells.Item(6, 8).Value = 141741

# --- Row 56: Barein ---
$ws.Cells.Item(56, 2).Value = 69361
$ws.Cells.Item(56, 3).Value = 586
$ws.Cells.Item(56, 4).Value = 62887
$ws.Cells.Item(56, 5).Value = 6232

# --- Row 84: Bulgaria ---
$ws.Cells.Item(84, 2).Value = 20055
$ws.Cells.Item(84, 3).Value = 58
$ws.Cells.Item(84, 4).Value = 14176
$ws.Cells.Item(84, 5).Value = 5083
$ws.Cells.Item(84, 7).Value = 7
$ws.Cells.Item(84, 8).Value = 796

# --- Row 85: Costa de Marfil ---
$ws.Cells.Item(85, 2).Value = 19629
$ws.Cells.Item(85, 3).Value = 29
$ws.Cells.Item(85, 4).Value = 19163
$ws.Cells.Item(85, 5).Value = 346

# --- Row 90: Tunez ---
$ws.Cells.Item(90, 2).Value = 16114
$ws.Cells.Item(90, 3).Value = 1722
$ws.Cells.Item(90, 5).Value = 10868
$ws.Cells.Item(90, 7).Value = 23
$ws.Cells.Item(90, 8).Value = 214

# --- Row 147: Guyana ---
$ws.Cells.Item(147, 2).Value = 2772
$ws.Cells.Item(147, 3).Value = 47
$ws.Cells.Item(147, 4).Value = 1564
$ws.Cells.Item(147, 5).Value = 1132
$ws.Cells.Item(147, 7).Value = 2
$ws.Cells.Item(147, 8).Value = 76

# --- Row 153: Yemen ---
$ws.Cells.Item(153, 4).Value = 1266
$ws.Cells.Item(153, 5).Value = 177
